$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.01477537992330724
$ws.Range("D2").Value = 0.01415548937485056
$ws.Range("E2").Value = 0.4214831381992497
$ws.Range("F2").Value = 0.5645336390059299
$ws.Range("G2").Value = 0.002391557367286488
$ws.Range("I2").Value = 0.4185668184037183
$ws.Range("K2").Value = 1.146770383261071
$ws.Range("N2").Value = 0.9335883308125972
$ws.Range("O2").Value = 1.841273413678095

$ws.Range("C3").Value = 0.01290205794836652
$ws.Range("D3").Value = 0.0129809571917221
$ws.Range("E3").Value = 0.3676435516299676
$ws.Range("F3").Value = 0.5564203375820469
$ws.Range("G3").Value = 0.002394573994975358
$ws.Range("I3").Value = 0.4132552972077264
$ws.Range("K3").Value = 1.004072333210246
$ws.Range("N3").Value = 0.9337865576973599
$ws.Range("O3").Value = 1.828548865374387

$ws.Range("C4").Value = 0.01174733087957236
$ws.Range("D4").Value = 0.01225527627693523
$ws.Range("E4").Value = 0.3346800167375221
$ws.Range("F4").Value = 0.5519166084740306
$ws.Range("G4").Value = 0.002396523281744133
$ws.Range("I4").Value = 0.4103525423954579
$ws.Range("K4").Value = 0.916272061756473
$ws.Range("N4").Value = 0.9342993659346988
$ws.Range("O4").Value = 1.822326892161357

$ws.Range("C5").Value = 0.01127566297977012
$ws.Range("D5").Value = 0.01195844569040361
$ws.Range("E5").Value = 0.3212689146550218
$ws.Range("F5").Value = 0.5502010188825253
$ws.Range("G5").Value = 0.002397342118314751
$ws.Range("I5").Value = 0.4092594540139203
$ws.Range("K5").Value = 0.8804482832820781
$ws.Range("N5").Value = 0.9346069210693244
$ws.Range("O5").Value = 1.820189829416194

$ws.Range("C6").Value = 0.01119727672614346
$ws.Range("D6").Value = 0.01190909091859638
$ws.Range("E6").Value = 0.3190432737337403
$ws.Range("F6").Value = 0.5499233640554664
$ws.Range("G6").Value = 0.002397479566659727
$ws.Range("I6").Value = 0.4090833614695697
$ws.Range("K6").Value = 0.8744971341569112
$ws.Range("N6").Value = 0.9346639515851365
$ws.Range("O6").Value = 1.819858989521947

$ws.Range("C7").Value = 0.01174097425258225
$ws.Range("D7").Value = 0.01225127757560784
$ws.Range("E7").Value = 0.3344990640945298
$ws.Range("F7").Value = 0.5518929873007536
$ws.Range("G7").Value = 0.002396534225542965
$ws.Range("I7").Value = 0.4103374374608606
$ws.Range("K7").Value = 0.9157891076374369
$ws.Range("N7").Value = 0.9343031142614322
$ws.Range("O7").Value = 1.822296459860809

$ws.Range("C8").Value = 0.01413040640738217
$ws.Range("D8").Value = 0.01375146213756295
$ws.Range("E8").Value = 0.4028985512062206
$ws.Range("F8").Value = 0.5616367164452001
$ws.Range("G8").Value = 0.002392577401919863
$ws.Range("I8").Value = 0.416660776803738
$ws.Range("K8").Value = 1.097607000462801
$ws.Range("N8").Value = 0.9335756451534962
$ws.Range("O8").Value = 1.836554804035359

$ws.Range("C9").Value = 0.01877952850743725
$ws.Range("D9").Value = 0.01665654239830872
$ws.Range("E9").Value = 0.5378681503396194
$ws.Range("F9").Value = 0.5845579591249361
$ws.Range("G9").Value = 0.002385584665110298
$ws.Range("I9").Value = 0.4319228239051256
$ws.Range("K9").Value = 1.452646574863024
$ws.Range("N9").Value = 0.9352431478781682
$ws.Range("O9").Value = 1.877216332760668

$ws.Range("C10").Value = 0.02217215385215354
$ws.Range("D10").Value = 0.0187673965114783
$ws.Range("E10").Value = 0.6376766563314362
$ws.Range("F10").Value = 0.6037566128758698
$ws.Range("G10").Value = 0.0023809093944
$ws.Range("I10").Value = 0.4449062604746743
$ws.Range("K10").Value = 1.712532936080663
$ws.Range("N10").Value = 0.9383441919515434
$ws.Range("O10").Value = 1.914947086865823

$ws.Range("C11").Value = 0.02371038538744585
$ws.Range("D11").Value = 0.01972235998799476
$ws.Range("E11").Value = 0.6832523784450757
$ws.Range("F11").Value = 0.6130100714400868
$ws.Range("G11").Value = 0.00237888180074074
$ws.Range("I11").Value = 0.4512028520698763
$ws.Range("K11").Value = 1.830545742776565
$ws.Range("N11").Value = 0.9401603653511614
$ws.Range("O11").Value = 1.933842830766849

$ws.Range("C12").Value = 0.02429212178347484
$ws.Range("D12").Value = 0.02008319934166281
$ws.Range("E12").Value = 0.700537756248707
$ws.Range("F12").Value = 0.6165894204194586
$ws.Range("G12").Value = 0.002378128188759909
$ws.Range("I12").Value = 0.4536437631232175
$ws.Range("K12").Value = 1.875202610670215
$ws.Range("N12").Value = 0.9409062413005103
$ws.Range("O12").Value = 1.941249106873869

$ws.Range("C13").Value = 0.02416686859636741
$ws.Range("D13").Value = 0.02000552143818624
$ws.Range("E13").Value = 0.6968138140740194
$ws.Range("F13").Value = 0.6158151880794946
$ws.Range("G13").Value = 0.002378289862513608
$ws.Range("I13").Value = 0.4531155488270713
$ws.Range("K13").Value = 1.865586406318755
$ws.Range("N13").Value = 0.9407430209556082
$ws.Range("O13").Value = 1.939642846106665

$ws.Range("C14").Value = 0.02375826054350227
$ws.Range("D14").Value = 0.01975206233169757
$ws.Range("E14").Value = 0.6846739070683725
$ws.Range("F14").Value = 0.6133030353125406
$ws.Range("G14").Value = 0.00237881951679215
$ws.Range("I14").Value = 0.4514025321081903
$ws.Range("K14").Value = 1.83422034164721
$ws.Range("N14").Value = 0.9402205649869302
$ws.Range("O14").Value = 1.934447111266593

$ws.Range("C15").Value = 0.0235078767495196
$ws.Range("D15").Value = 0.01959670836441063
$ws.Range("E15").Value = 0.6772414204940276
$ws.Range("F15").Value = 0.6117740859533427
$ws.Range("G15").Value = 0.002379145791260533
$ws.Range("I15").Value = 0.4503606333322736
$ws.Range("K15").Value = 1.815003500882028
$ws.Range("N15").Value = 0.9399081107642786
$ws.Range("O15").Value = 1.931297298593563

$ws.Range("C16").Value = 0.02207152210004892
$ws.Range("D16").Value = 0.018704879025222
$ws.Range("E16").Value = 0.6347018190504201
$ws.Range("F16").Value = 0.6031623823148635
$ws.Range("G16").Value = 0.002381043892370316
$ws.Range("I16").Value = 0.4445026521996098
$ws.Range("K16").Value = 1.704816133421332
$ws.Range("N16").Value = 0.938233644813522
$ws.Range("O16").Value = 1.913747210937629

$ws.Range("C17").Value = 0.02118904231310381
$ws.Range("D17").Value = 0.0181564009838624
$ws.Range("E17").Value = 0.6086507484823329
$ws.Range("F17").Value = 0.598012897136698
$ws.Range("G17").Value = 0.002382233674281094
$ws.Range("I17").Value = 0.4410092319229193
$ws.Range("K17").Value = 1.637164541122786
$ws.Range("N17").Value = 0.9373101406380897
$ws.Range("O17").Value = 1.903425666063299

$ws.Range("C18").Value = 0.02068098568845045
$ws.Range("D18").Value = 0.01784043598767227
$ws.Range("E18").Value = 0.5936829580297882
$ws.Range("F18").Value = 0.5950999540970798
$ws.Range("G18").Value = 0.002382927347709394
$ws.Range("I18").Value = 0.4390366278964777
$ws.Range("K18").Value = 1.598233393964222
$ws.Range("N18").Value = 0.9368171425011838
$ws.Range("O18").Value = 1.897651884816185

$ws.Range("C19").Value = 0.02050888522531125
$ws.Range("D19").Value = 0.01773337146386922
$ws.Range("E19").Value = 0.5886178234939194
$ws.Range("F19").Value = 0.5941220657083761
$ws.Range("G19").Value = 0.002383163820450066
$ws.Range("I19").Value = 0.4383750325710452
$ws.Range("K19").Value = 1.585048638276533
$ws.Range("N19").Value = 0.936656784343711
$ws.Range("O19").Value = 1.895724903971285

$ws.Range("C20").Value = 0.02128303341431348
$ws.Range("D20").Value = 0.01821483882124397
$ws.Range("E20").Value = 0.6114222486797587
$ws.Range("F20").Value = 0.5985560035248056
$ws.Range("G20").Value = 0.002382106053660506
$ws.Range("I20").Value = 0.4413773093462083
$ws.Range("K20").Value = 1.644368224109428
$ws.Range("N20").Value = 0.9374044993723629
$ws.Range("O20").Value = 1.904507539506199

$ws.Range("C21").Value = 0.02387829942901476
$ws.Range("D21").Value = 0.0198265309149761
$ws.Range("E21").Value = 0.6882389470510901
$ws.Range("F21").Value = 0.6140388687703364
$ws.Range("G21").Value = 0.00237866355972001
$ws.Range("I21").Value = 0.4519041493197236
$ws.Range("K21").Value = 1.843434197735746
$ws.Range("N21").Value = 0.9403724467071726
$ws.Range("O21").Value = 1.935966401405011

$ws.Range("C22").Value = 0.02557002248653362
$ws.Range("D22").Value = 0.02087527977514014
$ws.Range("E22").Value = 0.7386005215510778
$ws.Range("F22").Value = 0.624596754278798
$ws.Range("G22").Value = 0.002376496388728041
$ws.Range("I22").Value = 0.4591137234734006
$ws.Range("K22").Value = 1.973347895630809
$ws.Range("N22").Value = 0.942650923949401
$ws.Range("O22").Value = 1.957989693282684

$ws.Range("C23").Value = 0.02466753240764774
$ws.Range("D23").Value = 0.02031597085627368
$ws.Range("E23").Value = 0.7117065109468257
$ws.Range("F23").Value = 0.6189214818977433
$ws.Range("G23").Value = 0.002377645505827519
$ws.Range("I23").Value = 0.4552355417696319
$ws.Range("K23").Value = 1.904028219683141
$ws.Range("N23").Value = 0.9414039188168175
$ws.Range("O23").Value = 1.946100965919698

$ws.Range("C24").Value = 0.02124054223923366
$ws.Range("D24").Value = 0.01818842105848972
$ws.Range("E24").Value = 0.610169224616584
$ws.Range("F24").Value = 0.5983103169995303
$ws.Range("G24").Value = 0.002382163720915243
$ws.Range("I24").Value = 0.4412107900238738
$ws.Range("K24").Value = 1.64111155523608
$ws.Range("N24").Value = 0.937361721597
$ws.Range("O24").Value = 1.904017925587823

$ws.Range("C25").Value = 0.0175258141293213
$ws.Range("D25").Value = 0.01587469203952452
$ws.Range("E25").Value = 0.5012516340736823
$ws.Range("F25").Value = 0.5779454137143745
$ws.Range("G25").Value = 0.002387394839848831
$ws.Range("I25").Value = 0.4274849926994477
$ws.Range("K25").Value = 1.356764648001047
$ws.Range("N25").Value = 0.9344619011276194
$ws.Range("O25").Value = 1.864844845256897
